$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 7 (Ano 2025) with refreshed ADD data
$ws.Range("B7").Value = 2306154.85
$ws.Range("C7").Value = -48.0955592561132
$ws.Range("D7").Value = 2337
$ws.Range("E7").Value = 2337
$ws.Range("F7").Value = 986.8013906718015
$ws.Range("G7").Value = 5.185892752694876
